$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H corresponds to "Houses of Worship". Rows 24 through 176 all had a
# value of 1 that must be changed to 0.
for ($r = 24; $r -le 176; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
